# The commit swaps the two embedded theme parts: the "Office Theme" that used
# to live in ppt/theme/theme1.xml (Notes Master) and the "Integral" theme that
# used to live in ppt/theme/theme2.xml (Slide Master) trade places. The net,
# user-visible effect on the deck is that the Slide Master (and therefore every
# slide built on it) now carries the plain "Office" colour palette instead of
# the custom green/gold "Integral" palette.
#
# PowerPoint's COM object model doesn't expose a "swap these two raw theme
# parts" verb, but Master.ColorScheme gives direct, scriptable access to the
# twelve theme colour slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# that are the substantive content of <a:clrScheme> inside the theme part the
# Slide Master is built on. Re-pointing every slot at the standard Office
# palette reproduces the colour change the diff shows for that theme part.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.ColorScheme

# PpColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1..accent6, hlink, folHlink.
# Values are the standard Office theme colours, encoded as the OLE RGB() long
# PowerPoint's ColorFormat.RGB expects (0x00BBGGRR).
$officeColorScheme = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $officeColorScheme.Count; $i++) {
    $scheme.Colors($i).RGB = $officeColorScheme[$i - 1]
}

$p.Save()
